$d = $word.ActiveDocument

$replacements = @(
    @{old = "42×17=714";   new = "22×78=1716"},
    @{old = "38×90=3420";  new = "36×47=1692"},
    @{old = "22×30=660";   new = "44×15=660"},
    @{old = "35×73=2555";  new = "44×70=3080"},
    @{old = "95×46=4370";  new = "29×36=1044"},
    @{old = "19×57=1083";  new = "27×38=1026"},
    @{old = "73×16=1168";  new = "65×83=5395"},
    @{old = "24×88=2112";  new = "27×32=864"},
    @{old = "14×67=938";   new = "29×26=754"},
    @{old = "84×82=6888";  new = "43×94=4042"},
    @{old = "25×87=2175";  new = "55×50=2750"},
    @{old = "49×64=3136";  new = "29×94=2726"},
    @{old = "53×80=4240";  new = "68×96=6528"},
    @{old = "83×73=6059";  new = "14×89=1246"},
    @{old = "54×98=5292";  new = "95×79=7505"},
    @{old = "25×35=875";   new = "31×21=651"},
    @{old = "52×32=1664";  new = "85×28=2380"},
    @{old = "91×37=3367";  new = "67×94=6298"},
    @{old = "12×17=204";   new = "61×17=1037"},
    @{old = "87×31=2697";  new = "54×99=5346"},
    @{old = "66×48=3168";  new = "21×13=273"},
    @{old = "96×19=1824";  new = "12×56=672"},
    @{old = "54×17=918";   new = "90×99=8910"},
    @{old = "71×73=5183";  new = "74×20=1480"},
    @{old = "11×28=308";   new = "37×12=444"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
